# "Added Icons and Started Poster"
# Fill in actual time spent ("Temps réel (j)", column D) for the first
# few tasks, and move the active selection to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D4").Value = 0.5
$ws.Range("D5").Value = 0.5
$ws.Range("D6").Value = 0.25

$ws.Range("D7").Select()
